# RTM.xlsx - "update for next baseline"
# Add newly-authored Test Case references (Register / View Account) into the
# "Test Cases" column (E) of the Requirements Traceability Matrix, and trim a
# stale line from the Logout test-case list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Register feature (rows 2-4) ------------------------------------------------
$ws.Range("E2").Value = "TC_Register_001"
$ws.Range("E3").Value = "TC_Register_001"

$registerList = @(
    "TC_Register_002","TC_Register_003","TC_Register_004","TC_Register_005",
    "TC_Register_006","TC_Register_007","TC_Register_008","TC_Register_009",
    "TC_Register_010","TC_Register_011","TC_Register_012","TC_Register_013",
    "TC_Register_014","TC_Register_015","TC_Register_016","TC_Register_017",
    "TC_Register_018","TC_Register_019","TC_Register_020","TC_Register_021",
    "TC_Register_022","TC_Register_023","TC_Register_024"
) -join "`n"
$ws.Range("E4").Value = $registerList
$ws.Range("E4").WrapText = $true

# --- View Account feature (rows 10-13) ------------------------------------------
$viw1 = @("TC_ViwAccount_001","TC_ViwAccount_003","TC_ViwAccount_004","TC_ViwAccount_005","TC_ViwAccount_006") -join "`n"
$viw2 = @("TC_ViwAccount_009","TC_ViwAccount_011","TC_ViwAccount_012","TC_ViwAccount_013") -join "`n"
$ws.Range("E10").Value = $viw1
$ws.Range("E11").Value = $viw2
$ws.Range("E12").Value = "TC_ViwAccount_010"
$ws.Range("E13").Value = "TC_ViwAccount_002"

# --- Logout feature (row 32): drop the stale "TC_Logout _002" line --------------
$ws.Range("E32").Value = "TC_Logout _001_client`nTC_Logout _001_Admin`n"
$ws.Range("E32").Font.Bold = $true

# --- Row heights grew to accommodate the newly wrapped text ---------------------
$ws.Rows.Item(4).RowHeight = 409.6
$ws.Rows.Item(10).RowHeight = 105.75

# --- Restore the view: scrolled back to the top, final selection at E33 ---------
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("E33").Select() | Out-Null
